# Daily attendance processing - 2025-10-13 22:26:44
# For each data row, the "Recorded By" column (G) may contain a
# comma-separated list of recorders (e.g. "System, someone@example.com").
# Rotate the list so the last entry moves to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text -notmatch ",") {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Length -le 1) {
        continue
    }

    $lastPart = $parts[$parts.Length - 1]
    $rest = $parts[0..($parts.Length - 2)]
    $newParts = @($lastPart) + $rest
    $newText = $newParts -join ", "

    $cell.Value = $newText
}
